# Scheduled runner update: refresh Universalis market-price derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the Leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 38: Just Give Him a Serum (Hi-Potion of Strength)
$ws.Range("H38").Value = 361.69232
$ws.Range("I38").Value = 440.2
$ws.Range("J38").Value = 100
$ws.Range("K38").Value = 1320.6
$ws.Range("L38").Value = 300
$ws.Range("M38").Value = -948.5999999999999
$ws.Range("N38").Value = -1044

# Row 41: The Write Stuff (Enchanted Mythril Ink)
$ws.Range("H41").Value = 761.3333
$ws.Range("I41").Value = 875.9091
$ws.Range("J41").Value = 446.25
$ws.Range("K41").Value = 875.9091
$ws.Range("L41").Value = 446.25
$ws.Range("M41").Value = -435.9091
$ws.Range("N41").Value = -1326.25

# Row 47: Open Your Grimoire to Page 42 (Embossed Book of Silver)
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()

# Row 54: Arcane Arts for Dummies (Book of Mythril)
$ws.Range("H54").Value = 20000
$ws.Range("I54").Value = 20000
$ws.Range("K54").Value = 20000
$ws.Range("M54").Value = -19514

# Row 57: Quit Your Jib-jab (Gold Needle)
$ws.Range("H57").Value = 47455.4
$ws.Range("J57").Value = 47455.4
$ws.Range("L57").Value = 142366.2
$ws.Range("N57").Value = -143364.2

# Row 96: Scroll Down (Grade 1 Reisui of Intelligence)
$ws.Range("H96").Value = 669.8182
$ws.Range("I96").Value = 659.94446
$ws.Range("J96").Value = 714.25
$ws.Range("K96").Value = 1979.83338
$ws.Range("L96").Value = 2142.75
$ws.Range("M96").Value = -606.83338
$ws.Range("N96").Value = -4888.75

# Row 101: Edge of the Arcane (Cunning Craftsman's Tea)
$ws.Range("H101").Value = 25550796
$ws.Range("I101").Value = 770237.4
$ws.Range("J101").Value = 71571830
$ws.Range("K101").Value = 2310712.2
$ws.Range("L101").Value = 214715490
$ws.Range("M101").Value = -2309090.2
$ws.Range("N101").Value = -214718734

$ws = $wb.Worksheets.Item("ARM")
# Row 37: Get Shirty (Steel Chainmail)
$ws.Range("H37").Value = 4189956
$ws.Range("J37").Value = 24947.1
$ws.Range("L37").Value = 24947.1
$ws.Range("N37").Value = -25493.1

# Row 45: Hollow Hallmarks (Mythril Ingot)
$ws.Range("H45").Value = 13755.7
$ws.Range("I45").Value = 36169
$ws.Range("K45").Value = 36169
$ws.Range("M45").Value = -35792

# Row 110: Scheduled Maintenance (Deepgold Ingot)
$ws.Range("H110").Value = 1092.5714
$ws.Range("I110").Value = 1241.3334
$ws.Range("K110").Value = 1241.3334
$ws.Range("M110").Value = 803.6666

$ws = $wb.Worksheets.Item("BSM")
# Row 2: Proly Hatchet (Bronze Hatchet)
$ws.Range("H2").Value = 29995
$ws.Range("J2").Value = 29995
$ws.Range("L2").Value = 29995
$ws.Range("N2").Value = -30221

$ws = $wb.Worksheets.Item("CRP")
# Row 19: Shielding Sales (Square Ash Shield)
$ws.Range("H19").Value = 1343
$ws.Range("I19").Value = 576.25
$ws.Range("J19").Value = 1956.4
$ws.Range("K19").Value = 576.25
$ws.Range("L19").Value = 1956.4
$ws.Range("M19").Value = -406.25
$ws.Range("N19").Value = -2296.4

# Row 24: What You Need (Square Ash Shield)
$ws.Range("H24").Value = 1343
$ws.Range("I24").Value = 576.25
$ws.Range("J24").Value = 1956.4
$ws.Range("K24").Value = 576.25
$ws.Range("L24").Value = 1956.4
$ws.Range("M24").Value = -406.25
$ws.Range("N24").Value = -2296.4

# Row 31: Wall Not Found (Walnut Lumber)
$ws.Range("H31").Value = 37733.844
$ws.Range("I31").Value = 33945.156
$ws.Range("K31").Value = 33945.156
$ws.Range("M31").Value = -33650.156

# Row 34: Armoires of the Rich and Famous (Walnut Lumber)
$ws.Range("H34").Value = 37733.844
$ws.Range("I34").Value = 33945.156
$ws.Range("K34").Value = 33945.156
$ws.Range("M34").Value = -33743.156

# Row 52: Spin It Like You Mean It (Mahogany Spinning Wheel)
$ws.Range("H52").Value = 45935.5
$ws.Range("I52").Value = 41873.5
$ws.Range("J52").Value = 49997.5
$ws.Range("K52").Value = 41873.5
$ws.Range("L52").Value = 49997.5
$ws.Range("M52").Value = -41579.5
$ws.Range("N52").Value = -50585.5

# Row 58: You Do the Heavy Lifting (Mahogany Lumber)
$ws.Range("H58").Value = 910.36365
$ws.Range("I58").Value = 900
$ws.Range("J58").Value = 1014
$ws.Range("K58").Value = 900
$ws.Range("L58").Value = 1014
$ws.Range("M58").Value = -697
$ws.Range("N58").Value = -1420

# Row 107: Built to Last (White Oak Lumber)
$ws.Range("H107").Value = 2077.7144
$ws.Range("I107").Value = 1864
$ws.Range("K107").Value = 1864
$ws.Range("M107").Value = 56

# Row 136: Turali Quality (Dark Mahogany Lumber)
$ws.Range("H136").Value = 910.36365
$ws.Range("I136").Value = 900
$ws.Range("J136").Value = 1014
$ws.Range("K136").Value = 2700
$ws.Range("L136").Value = 3042
$ws.Range("M136").Value = -150
$ws.Range("N136").Value = -8142

$ws = $wb.Worksheets.Item("CUL")
# Row 70: Persona non Gratin (Dhalmel Gratin)
$ws.Range("H70").Value = 3501.4285
$ws.Range("I70").Value = 756
$ws.Range("J70").Value = 4599.6
$ws.Range("K70").Value = 2268
$ws.Range("L70").Value = 13798.8
$ws.Range("M70").Value = -1953
$ws.Range("N70").Value = -14428.8

# Row 73: Recipe for Disaster (L) (Dhalmel Gratin)
$ws.Range("H73").Value = 3501.4285
$ws.Range("I73").Value = 756
$ws.Range("J73").Value = 4599.6
$ws.Range("K73").Value = 2268
$ws.Range("L73").Value = 13798.8
$ws.Range("M73").Value = -1176
$ws.Range("N73").Value = -15982.8

# Row 80: Saucy for a Suitor (Hollandaise Sauce)
$ws.Range("H80").Value = 2400
$ws.Range("J80").Value = 2400
$ws.Range("L80").Value = 7200
$ws.Range("N80").Value = -9072

# Row 83: Saved by the Sauce (L) (Hollandaise Sauce)
$ws.Range("H83").Value = 2400
$ws.Range("J83").Value = 2400
$ws.Range("L83").Value = 21600
$ws.Range("N83").Value = -30960

# Row 88: Don't Let It Fall Apart (Liver-cheese Sandwich)
$ws.Range("H88").Value = 10332.4
$ws.Range("J88").Value = 10332.4
$ws.Range("L88").Value = 30997.2
$ws.Range("N88").Value = -31853.2

# Row 91: Better Come Back with a Sandwich (L) (Liver-cheese Sandwich)
$ws.Range("H91").Value = 10332.4
$ws.Range("J91").Value = 10332.4
$ws.Range("L91").Value = 30997.2
$ws.Range("N91").Value = -33961.2

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell (Hardsilver Ingot)
$ws.Range("H80").Value = 3879.842
$ws.Range("I80").Value = 3181.2727
$ws.Range("J80").Value = 4840.375
$ws.Range("K80").Value = 3181.2727
$ws.Range("L80").Value = 4840.375
$ws.Range("M80").Value = -2183.2727
$ws.Range("N80").Value = -6836.375

# Row 83: With a Noise That Reaches Heaven (L) (Hardsilver Ingot)
$ws.Range("H83").Value = 3879.842
$ws.Range("I83").Value = 3181.2727
$ws.Range("J83").Value = 4840.375
$ws.Range("K83").Value = 15906.3635
$ws.Range("L83").Value = 24201.875
$ws.Range("M83").Value = -10914.3635
$ws.Range("N83").Value = -34185.875

# Row 95: Chain of Command (Koppranickel Temple Chain)
$ws.Range("H95").Value = 28836
$ws.Range("J95").Value = 28836
$ws.Range("L95").Value = 28836
$ws.Range("N95").Value = -34328

# Row 122: Awarding Academic Excellence (Ametrine)
$ws.Range("H122").Value = 1330.2
$ws.Range("I122").Value = 1330.2
$ws.Range("K122").Value = 3990.6
$ws.Range("M122").Value = -1540.6

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban (Leather)
$ws.Range("H7").Value = 25557.666
$ws.Range("I7").Value = 27951.23
$ws.Range("K7").Value = 27951.23
$ws.Range("M7").Value = -27839.23

# Row 22: Skin off Their Backs (Aldgoat Leather)
$ws.Range("H22").Value = 1079.5
$ws.Range("I22").Value = 1049
$ws.Range("K22").Value = 1049
$ws.Range("M22").Value = -754

# Row 27: Fire and Hide (Aldgoat Leather)
$ws.Range("H27").Value = 1079.5
$ws.Range("I27").Value = 1049
$ws.Range("K27").Value = 1049
$ws.Range("M27").Value = -942

# Row 40: Best Served Toad (Toad Leather)
$ws.Range("H40").Value = 8649.308000000001
$ws.Range("I40").Value = 8864.434999999999
$ws.Range("K40").Value = 8864.434999999999
$ws.Range("M40").Value = -8728.434999999999

# Row 46: Supply Side Logic (Boar Leather)
$ws.Range("H46").Value = 24478.475
$ws.Range("I46").Value = 61843.43
$ws.Range("J46").Value = 2682.25
$ws.Range("K46").Value = 61843.43
$ws.Range("L46").Value = 2682.25
$ws.Range("M46").Value = -61655.43
$ws.Range("N46").Value = -3058.25

# Row 55: It's Not a Job, It's a Calling (Peiste Leather)
$ws.Range("H55").Value = 328.92
$ws.Range("J55").Value = 265.42856
$ws.Range("L55").Value = 265.42856
$ws.Range("N55").Value = -611.4285600000001

# Row 93: Hide to Go Seek (Gagana Leather)
$ws.Range("H93").Value = 9450.65
$ws.Range("I93").Value = 1168.921
$ws.Range("K93").Value = 1168.921
$ws.Range("M93").Value = 79.07899999999995

# Row 122: Hell on Leather (Gaja Leather)
$ws.Range("H122").Value = 147578.72
$ws.Range("I122").Value = 204000.2
$ws.Range("K122").Value = 612000.6000000001
$ws.Range("M122").Value = -609550.6000000001

# Row 126: Battered Books (Saiga Leather)
$ws.Range("H126").Value = 25557.666
$ws.Range("I126").Value = 27951.23
$ws.Range("K126").Value = 83853.69
$ws.Range("M126").Value = -81383.69

# Row 132: Tenets of Tanning (Silver Lobo Leather)
$ws.Range("H132").Value = 3493.5278
$ws.Range("I132").Value = 3376.7585
$ws.Range("J132").Value = 3977.2856
$ws.Range("K132").Value = 10130.2755
$ws.Range("L132").Value = 11931.8568
$ws.Range("M132").Value = -7600.2755
$ws.Range("N132").Value = -16991.8568

$ws = $wb.Worksheets.Item("WVR")
# Row 39: By the Short Hairs (Velveteen Robe)
$ws.Range("H39").Value = 21499
$ws.Range("J39").Value = 21499
$ws.Range("L39").Value = 21499
$ws.Range("N39").Value = -22325

# Row 126: A Polished Purchase (Snow Linen)
$ws.Range("H126").Value = 1864.5294
$ws.Range("I126").Value = 1826.5333
$ws.Range("J126").Value = 2149.5
$ws.Range("K126").Value = 5479.5999
$ws.Range("M126").Value = -3009.5999
$ws.Range("N126").Value = -11388.5
